$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was added for "Cuatro cascos verde" pimiento at
# Macroferia Regional de Talca. The source sheet keeps rows ordered, so the
# new record is inserted at row 362 and every existing record from the old
# row 362 through row 380 shifts down by one row (to rows 363-381).
$ws.Rows.Item(362).Insert()

# Fill in the values for the newly inserted row 362.
$ws.Range("A362").Value = 5
$ws.Range("B362").Value = "Macroferia Regional de Talca"
$ws.Range("C362").Value = "Maule"
$ws.Range("D362").Value = 44516
$ws.Range("E362").Value = 7
$ws.Range("F362").Value = 100112002
$ws.Range("G362").Value = "Pimiento"
$ws.Range("H362").Value = "Cuatro cascos verde"
$ws.Range("I362").Value = "Primera"
$ws.Range("J362").Value = 200
$ws.Range("K362").Value = 18000
$ws.Range("L362").Value = 18000
$ws.Range("M362").Value = 18000
$ws.Range("N362").Value = "`$/caja 15 kilos"
$ws.Range("O362").Value = "Región del Maule"
$ws.Range("P362").Value = 1200
$ws.Range("Q362").Value = 15
$ws.Range("R362").Value = "Hortaliza"
